# Update "想去人数" (interested-count) figures that changed between crawls.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1104
$ws1.Range("F8").Value  = 823
$ws1.Range("F10").Value = 6006
$ws1.Range("F14").Value = 5794
$ws1.Range("F15").Value = 5794
$ws1.Range("F20").Value = 1621
$ws1.Range("F24").Value = 1308
$ws1.Range("F31").Value = 3846

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2209
$ws3.Range("F5").Value = 149

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2209
$ws4.Range("F7").Value  = 1104
$ws4.Range("F12").Value = 823
$ws4.Range("F13").Value = 149
$ws4.Range("F15").Value = 6006
$ws4.Range("F24").Value = 5794
$ws4.Range("F25").Value = 5794
$ws4.Range("F30").Value = 1621
$ws4.Range("F34").Value = 1308
$ws4.Range("F46").Value = 3846
